$p = $ppt.ActivePresentation
Write-Output ("Presentations.Count=" + $ppt.Presentations.Count)
try {
    $p2 = $ppt.Presentations.Open($p.FullName)
    Write-Output "open ok"
    Write-Output ("Presentations.Count now=" + $ppt.Presentations.Count)
} catch {
    Write-Output ("open err: " + $_.Exception.Message)
}
